$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D can look like numbers (e.g. "1.00", "0.696") to Excel's
# automatic type detection, but the source data keeps them as literal text.
# Force the cells to Text format before writing so COM stores them as strings,
# then restore the "Normal" cell style afterwards so no stray formatting remains
# on the cell itself (only the underlying styles table gains an unused Text
# number-format entry, which no cell references).
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D10", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D26", "D28", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D40", "D41", "D45", "D46", "D47", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "89.881.89"
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("D3").Value = "3.196.66"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "216.56"
$ws.Range("E5").Value = "  +5.07%  "
$ws.Range("D6").Value = "622.29"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("D8").Value = "0.696"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "3.194.75"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  +4.72%  "
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "5.37"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "3.786.68"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "89.650.96"
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("D17").Value = "32.68"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "3.211.32"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").Value = "0.0000232"
$ws.Range("E19").Value = "  +74.50%  "
$ws.Range("D20").Value = "3.34"
$ws.Range("E20").Value = "  +12.38%  "
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").Value = "431.44"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").Value = "11.57"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").Value = "75.17"
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "0.156"
$ws.Range("E31").Value = "  -10.41%  "
$ws.Range("D32").Value = "4.03"
$ws.Range("E32").Value = "  +33.79%  "
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("D34").Value = "533.80"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").Value = "6.84"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("D37").Value = "1.27"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").Value = "22.16"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "0.126"
$ws.Range("E41").Value = "  -6.60%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("D45").Value = "150.53"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("D46").Value = "171.07"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("D47").Value = "43.16"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("E48").Value = "  -4.56%  "
$ws.Range("E49").Value = "  -4.28%  "
$ws.Range("D50").Value = "0.731"
$ws.Range("E50").Value = "  +3.55%  "
$ws.Range("E51").Value = "  +1.02%  "

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
